# Update "paises.xlsx" (sheet "Pais") with the new COVID-19 snapshot data:
#  - refresh the "last updated" timestamp
#  - update case counts for a number of countries
#  - fix the row order for three country pairs whose data/labels were
#    swapped (Austria/Emiratos Arabes Unidos, Belice/Nueva Caledonia,
#    Curazao/Dominica)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 13:04"
$ws.Range("B4").Value = 1263243
$ws.Range("C4").Value = 151
$ws.Range("E4").Value = 975325
$ws.Range("B5").Value = 256855
$ws.Range("C5").Value = 3173
$ws.Range("D5").Value = 163919
$ws.Range("E5").Value = 66866
$ws.Range("G5").Value = 213
$ws.Range("H5").Value = 26070
$ws.Range("B21").Value = 30126
$ws.Range("C21").Value = 66
$ws.Range("E21").Value = 2621
$ws.Range("B31").Value = 18890
$ws.Range("C31").Value = 918
$ws.Range("D31").Value = 2286
$ws.Range("E31").Value = 16592

# Austria / Emiratos Arabes Unidos swap rows 33-34
$ws.Range("A33").Value = "Emiratos Arabes Unidos"
$ws.Range("B33").Value = 16240
$ws.Range("C33").Value = 502
$ws.Range("D33").Value = 3572
$ws.Range("E33").Value = 12503
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 8
$ws.Range("H33").Value = 165
$ws.Range("A34").Value = "Austria"
$ws.Range("B34").Value = 15752
$ws.Range("C34").Value = 68
$ws.Range("D34").Value = 13698
$ws.Range("E34").Value = 1445
$ws.Range("F34").Value = 92
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 609
$ws.Range("B52").Value = 6896
$ws.Range("C52").Value = 21
$ws.Range("E52").Value = 759
$ws.Range("F52").Value = 23
$ws.Range("B55").Value = 5673
$ws.Range("C55").Value = 100
$ws.Range("E55").Value = 1921
$ws.Range("B56").Value = 5505
$ws.Range("C56").Value = 97
$ws.Range("D56").Value = 2124
$ws.Range("E56").Value = 3198
$ws.Range("B59").Value = 4509
$ws.Range("C59").Value = 87
$ws.Range("D59").Value = 1450
$ws.Range("E59").Value = 3029
$ws.Range("B63").Value = 3563
$ws.Range("C63").Value = 171
$ws.Range("D63").Value = 468
$ws.Range("E63").Value = 2989
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 106
$ws.Range("B100").Value = 842
$ws.Range("C100").Value = 10
$ws.Range("D100").Value = 605
$ws.Range("D174").Value = 40
$ws.Range("E174").Value = 5
$ws.Range("D188").Value = 9
$ws.Range("E188").Value = 10

# Belice / Nueva Caledonia swap rows 191-192
$ws.Range("A191").Value = "Nueva Caledonia"
$ws.Range("D191").Value = 18
$ws.Range("H191").Value = 0
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

# Curazao / Dominica swap rows 198-199
$ws.Range("A198").Value = "Dominica"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 0
$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 13
$ws.Range("H199").Value = 1
